$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename move id / move_id values to zero-padded two-digit form (move_1 -> move_01 etc.)
$ws.Range("A2").Value = "z0bug.move_01_1"
$ws.Range("B2").Value = "z0bug.move_01"
$ws.Range("A3").Value = "z0bug.move_01_2"
$ws.Range("B3").Value = "z0bug.move_01"
$ws.Range("A4").Value = "z0bug.move_02_1"
$ws.Range("B4").Value = "z0bug.move_02"
$ws.Range("A5").Value = "z0bug.move_02_2"
$ws.Range("B5").Value = "z0bug.move_02"
$ws.Range("A6").Value = "z0bug.move_03_1"
$ws.Range("B6").Value = "z0bug.move_03"
$ws.Range("A7").Value = "z0bug.move_03_2"
$ws.Range("B7").Value = "z0bug.move_03"
$ws.Range("A8").Value = "z0bug.move_03_3"
$ws.Range("B8").Value = "z0bug.move_03"
$ws.Range("A9").Value = "z0bug.move_04_1"
$ws.Range("B9").Value = "z0bug.move_04"
$ws.Range("A10").Value = "z0bug.move_04_2"
$ws.Range("B10").Value = "z0bug.move_04"
$ws.Range("A11").Value = "z0bug.move_05_1"
$ws.Range("B11").Value = "z0bug.move_05"
$ws.Range("A12").Value = "z0bug.move_05_2"
$ws.Range("B12").Value = "z0bug.move_05"
$ws.Range("A13").Value = "z0bug.move_06_1"
$ws.Range("B13").Value = "z0bug.move_06"
$ws.Range("A14").Value = "z0bug.move_06_2"
$ws.Range("B14").Value = "z0bug.move_06"
$ws.Range("A15").Value = "z0bug.move_07_1"
$ws.Range("B15").Value = "z0bug.move_07"
$ws.Range("A16").Value = "z0bug.move_07_2"
$ws.Range("B16").Value = "z0bug.move_07"
$ws.Range("A17").Value = "z0bug.move_07_3"
$ws.Range("B17").Value = "z0bug.move_07"
$ws.Range("A18").Value = "z0bug.move_08_1"
$ws.Range("B18").Value = "z0bug.move_08"
$ws.Range("A19").Value = "z0bug.move_08_2"
$ws.Range("B19").Value = "z0bug.move_08"
$ws.Range("A20").Value = "z0bug.move_09_1"
$ws.Range("B20").Value = "z0bug.move_09"
$ws.Range("A21").Value = "z0bug.move_09_2"
$ws.Range("B21").Value = "z0bug.move_09"

# Adjust column widths for columns A and B (slightly wider after edit)
$ws.Columns.Item(1).ColumnWidth = 16.67
$ws.Columns.Item(2).ColumnWidth = 13.83

# Move the active selection to C5
$ws.Range("C5").Select()
